$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 43 (hunk 0)
$ws.Range("H43").Value = 1432.3334
$ws.Range("I43").Value = 1313.1538
$ws.Range("J43").Value = 1626
$ws.Range("K43").Value = 1313.1538
$ws.Range("L43").Value = 1626
$ws.Range("M43").Value = -1244.1538
$ws.Range("N43").Value = -1764

# row 93 (hunk 1)
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# row 114 (hunk 2)
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# row 125 (hunk 3)
$ws.Range("H125").Value = 25106.4
$ws.Range("I125").Value = 25106.4
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 225957.6
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -223497.6
$ws.Range("N125").ClearContents()

# row 135 (hunk 4)
$ws.Range("H135").Value = 1055.6129
$ws.Range("I135").Value = 407.6
$ws.Range("K135").Value = 3668.4
$ws.Range("M135").Value = -1133.4

# row 141 (hunk 5)
$ws.Range("H141").Value = 4920.263
$ws.Range("I141").Value = 2038.0358
$ws.Range("J141").Value = 12990.5
$ws.Range("K141").Value = 6114.107400000001
$ws.Range("L141").Value = 38971.5
$ws.Range("M141").Value = -934.1074000000008
$ws.Range("N141").Value = -49331.5

$ws = $wb.Worksheets.Item("ARM")
# row 43 (hunk 6)
$ws.Range("H43").Value = 9377
$ws.Range("J43").Value = 9377
$ws.Range("L43").Value = 9377
$ws.Range("N43").Value = -10003

# row 44 (hunk 7)
$ws.Range("H44").Value = 24966.666
$ws.Range("J44").Value = 24966.666
$ws.Range("L44").Value = 24966.666
$ws.Range("N44").Value = -25942.666

# row 61 (hunk 8)
$ws.Range("H61").Value = 14927455
$ws.Range("I61").Value = 19609754
$ws.Range("J61").Value = 2627.875
$ws.Range("K61").Value = 19609754
$ws.Range("L61").Value = 2627.875
$ws.Range("M61").Value = -19609542
$ws.Range("N61").Value = -3051.875

# row 74 (hunk 9)
$ws.Range("H74").Value = 10061.167
$ws.Range("I74").Value = 1618.4
$ws.Range("J74").Value = 16091.714
$ws.Range("K74").Value = 1618.4
$ws.Range("L74").Value = 16091.714
$ws.Range("M74").Value = -744.4000000000001
$ws.Range("N74").Value = -17839.714

# row 77 (hunk 10)
$ws.Range("H77").Value = 10061.167
$ws.Range("I77").Value = 1618.4
$ws.Range("J77").Value = 16091.714
$ws.Range("K77").Value = 8092
$ws.Range("L77").Value = 80458.57000000001
$ws.Range("M77").Value = -3724
$ws.Range("N77").Value = -89194.57000000001

# row 136 (hunk 11)
$ws.Range("H136").Value = 14927455
$ws.Range("I136").Value = 19609754
$ws.Range("J136").Value = 2627.875
$ws.Range("K136").Value = 58829262
$ws.Range("L136").Value = 7883.625
$ws.Range("M136").Value = -58826712
$ws.Range("N136").Value = -12983.625

$ws = $wb.Worksheets.Item("CRP")
# row 10 (hunk 12)
$ws.Range("H10").Value = 309.75
$ws.Range("I10").Value = 309.75
$ws.Range("K10").Value = 309.75
$ws.Range("M10").Value = -170.75

# row 16 (hunk 13)
$ws.Range("H16").Value = 3147.8462
$ws.Range("I16").Value = 1003
$ws.Range("J16").Value = 5650.1665
$ws.Range("K16").Value = 1003
$ws.Range("L16").Value = 5650.1665
$ws.Range("M16").Value = -716
$ws.Range("N16").Value = -6224.1665

# row 19 (hunk 14)
$ws.Range("H19").Value = 392.85715
$ws.Range("I19").Value = 210
$ws.Range("J19").Value = 850
$ws.Range("K19").Value = 210
$ws.Range("L19").Value = 850
$ws.Range("M19").Value = -40
$ws.Range("N19").Value = -1190

# row 24 (hunk 15)
$ws.Range("H24").Value = 392.85715
$ws.Range("I24").Value = 210
$ws.Range("J24").Value = 850
$ws.Range("K24").Value = 210
$ws.Range("L24").Value = 850
$ws.Range("M24").Value = -40
$ws.Range("N24").Value = -1190

# row 113 (hunk 16)
$ws.Range("H113").Value = 3147.8462
$ws.Range("I113").Value = 1003
$ws.Range("J113").Value = 5650.1665
$ws.Range("K113").Value = 1003
$ws.Range("L113").Value = 5650.1665
$ws.Range("M113").Value = 1167
$ws.Range("N113").Value = -9990.166499999999

$ws = $wb.Worksheets.Item("GSM")
# row 3 (hunk 17)
$ws.Range("H3").Value = 318.33334
$ws.Range("I3").Value = 318.33334
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 318.33334
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -202.33334
$ws.Range("N3").ClearContents()

# row 7 (hunk 18)
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

# row 8 (hunk 19)
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

# row 9 (hunk 20)
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()

# row 10 (hunk 21)
$ws.Range("H10").Value = 668334.3
$ws.Range("J10").Value = 4000
$ws.Range("L10").Value = 4000
$ws.Range("N10").Value = -4338

# row 11 (hunk 22)
$ws.Range("H11").Value = 27895.25
$ws.Range("I11").Value = 900
$ws.Range("J11").Value = 44092.4
$ws.Range("K11").Value = 900
$ws.Range("L11").Value = 44092.4
$ws.Range("M11").Value = -761
$ws.Range("N11").Value = -44370.4

# row 12 (hunk 23)
$ws.Range("H12").Value = 433.33334
$ws.Range("J12").Value = 433.33334
$ws.Range("L12").Value = 433.33334
$ws.Range("N12").Value = -713.33334

# row 13 (hunk 24)
$ws.Range("H13").Value = 201.66667
$ws.Range("I13").Value = 201.66667
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 201.66667
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -62.66667000000001
$ws.Range("N13").ClearContents()

# row 80 (hunk 25)
$ws.Range("H80").Value = 2375.4443
$ws.Range("I80").Value = 1866.4
$ws.Range("J80").Value = 2674.8823
$ws.Range("K80").Value = 1866.4
$ws.Range("L80").Value = 2674.8823
$ws.Range("M80").Value = -868.4000000000001
$ws.Range("N80").Value = -4670.8823

# row 83 (hunk 26)
$ws.Range("H83").Value = 2375.4443
$ws.Range("I83").Value = 1866.4
$ws.Range("J83").Value = 2674.8823
$ws.Range("K83").Value = 9332
$ws.Range("L83").Value = 13374.4115
$ws.Range("M83").Value = -4340
$ws.Range("N83").Value = -23358.4115

# row 107 (hunk 27)
$ws.Range("H107").Value = 165.9375
$ws.Range("I107").Value = 158.72728
$ws.Range("J107").Value = 181.8
$ws.Range("K107").Value = 158.72728
$ws.Range("L107").Value = 181.8
$ws.Range("M107").Value = 1761.27272
$ws.Range("N107").Value = -4021.8

# row 132 (hunk 28)
$ws.Range("H132").Value = 8598.477000000001
$ws.Range("I132").Value = 10919.714
$ws.Range("K132").Value = 32759.142
$ws.Range("M132").Value = -30229.142

$ws = $wb.Worksheets.Item("LTW")
# row 41 (hunk 29)
$ws.Range("H41").Value = 7000
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

# row 136 (hunk 30)
$ws.Range("H136").Value = 5905.1304
$ws.Range("I136").Value = 1467.3889
$ws.Range("J136").Value = 21881
$ws.Range("K136").Value = 4402.1667
$ws.Range("L136").Value = 65643
$ws.Range("M136").Value = -1852.1667
$ws.Range("N136").Value = -70743

$ws = $wb.Worksheets.Item("WVR")
# row 136 (hunk 31)
$ws.Range("H136").Value = 1490.9697
$ws.Range("I136").Value = 800.7368
$ws.Range("J136").Value = 2427.7144
$ws.Range("K136").Value = 2402.2104
$ws.Range("L136").Value = 7283.1432
$ws.Range("M136").Value = 147.7896000000001
$ws.Range("N136").Value = -12383.1432
